$wb = $excel.ActiveWorkbook
$poland = $wb.Worksheets.Item("Poland")
$poland.Copy([System.Reflection.Missing]::Value, $poland)
$newSheet = $wb.Worksheets.Item($poland.Index + 1)
$newSheet.Name = "UK"

$newSheet.Range("B4").Value = "NGC-2741/T3343/T3339"
$newSheet.Range("B2").Value = "UK Market"

$newSheet.Range("B4").Select()

Write-Output $newSheet.Name
Write-Output $wb.Worksheets.Count
